$wb = $excel.ActiveWorkbook

# Sheet "展览" - update F2 287 -> 289 and F5 903 -> 907
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 289
$ws1.Range("F5").Value = 907

# Sheet "全部类型" - update F2 287 -> 289 and F5 903 -> 907
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 289
$ws4.Range("F5").Value = 907
